# Insert a new data row at row 377 (pushing the existing rows 377-445
# down to 378-446), then populate the new row with its values.
#
# The new row carries the same Mercado/Region/Categoria/Variedad/Calidad/
# Volumen/Clasificacion values as the row that used to be at 377, but with
# its own Fecha, Precio minimo/maximo/promedio, Unidad de comercializacion,
# Origen, Precio $/Kg and Kg o Unidades.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 377:445 down to 378:446, leaving row 377 free for the new record.
$ws.Rows.Item(377).Insert()

$ws.Cells.Item(377, 1).Value = 6
$ws.Cells.Item(377, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(377, 3).Value = "Metropolitana"
$ws.Cells.Item(377, 4).Value = 44711
$ws.Cells.Item(377, 5).Value = 13
$ws.Cells.Item(377, 6).Value = 100112043
$ws.Cells.Item(377, 7).Value = "Pepino ensalada"
$ws.Cells.Item(377, 8).Value = "Sin especificar"
$ws.Cells.Item(377, 9).Value = "Primera"
$ws.Cells.Item(377, 10).Value = 400
$ws.Cells.Item(377, 11).Value = 18000
$ws.Cells.Item(377, 12).Value = 19000
$ws.Cells.Item(377, 13).Value = 18425
$ws.Cells.Item(377, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(377, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(377, 16).Value = 307
$ws.Cells.Item(377, 17).Value = 60
$ws.Cells.Item(377, 18).Value = "Hortaliza"
